$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf2"
$ws.Range("C2").Value = "Gpc4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.005243333333333
$ws.Range("H2").Value = 3.01573
$ws.Range("I2").Value = 0.07224874268505826
$ws.Range("J2").Value = 0.07224874268505825
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.269799333333333
$ws.Range("N2").Value = 27.809398
$ws.Range("O2").Value = 0.1628450003986012
$ws.Range("P2").Value = 0.1628450003986012
$ws.Range("Q2").Value = 9.31840398117111
$ws.Range("R2").Value = 83.86563583053999
$ws.Range("S2").Value = 0.01176534653134675
$ws.Range("T2").Value = 0.01176534653134675

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf2"
$ws.Range("C3").Value = "Gpc4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.005243333333333
$ws.Range("H3").Value = 3.01573
$ws.Range("I3").Value = 0.07224874268505826
$ws.Range("J3").Value = 0.07224874268505825
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 37.71549866666666
$ws.Range("N3").Value = 113.146496
$ws.Range("O3").Value = 0.6625580742963342
$ws.Range("P3").Value = 0.6625580742963342
$ws.Range("Q3").Value = 37.91325359800888
$ws.Range("R3").Value = 341.2192823820799
$ws.Range("S3").Value = 0.04786898782374356
$ws.Range("T3").Value = 0.04786898782374355

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf2"
$ws.Range("C4").Value = "Gpc4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.005243333333333
$ws.Range("H4").Value = 3.01573
$ws.Range("I4").Value = 0.07224874268505826
$ws.Range("J4").Value = 0.07224874268505825
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04982966666666667
$ws.Range("N4").Value = 0.149489
$ws.Range("O4").Value = 0.0008753708463802955
$ws.Range("P4").Value = 0.0008753708463802955
$ws.Range("Q4").Value = 0.05009094021888889
$ws.Range("R4").Value = 0.45081846197
$ws.Range("S4").Value = 0.00006324444303413163
$ws.Range("T4").Value = 0.00006324444303413162

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf2"
$ws.Range("C5").Value = "Gpc4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.005243333333333
$ws.Range("H5").Value = 3.01573
$ws.Range("I5").Value = 0.07224874268505826
$ws.Range("J5").Value = 0.07224874268505825
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.888937
$ws.Range("N5").Value = 29.666811
$ws.Range("O5").Value = 0.1737215544586843
$ws.Range("P5").Value = 0.1737215544586843
$ws.Range("Q5").Value = 9.940787993003333
$ws.Range("R5").Value = 89.46709193702999
$ws.Range("S5").Value = 0.01255116388693382
$ws.Range("T5").Value = 0.01255116388693382

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf2"
$ws.Range("C6").Value = "Gpc4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.25983933333333
$ws.Range("H6").Value = 30.779518
$ws.Range("I6").Value = 0.7373940889775011
$ws.Range("J6").Value = 0.737394088977501
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.269799333333333
$ws.Range("N6").Value = 27.809398
$ws.Range("O6").Value = 0.1628450003986012
$ws.Range("P6").Value = 0.1628450003986012
$ws.Range("Q6").Value = 95.10665181224044
$ws.Range("R6").Value = 855.959866310164
$ws.Range("S6").Value = 0.1200809407134674
$ws.Range("T6").Value = 0.1200809407134673

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf2"
$ws.Range("C7").Value = "Gpc4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.25983933333333
$ws.Range("H7").Value = 30.779518
$ws.Range("I7").Value = 0.7373940889775011
$ws.Range("J7").Value = 0.737394088977501
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 37.71549866666666
$ws.Range("N7").Value = 113.146496
$ws.Range("O7").Value = 0.6625580742963342
$ws.Range("P7").Value = 0.6625580742963342
$ws.Range("Q7").Value = 386.9549566965475
$ws.Range("R7").Value = 3482.594610268928
$ws.Range("S7").Value = 0.4885664075904328
$ws.Range("T7").Value = 0.4885664075904328

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fgf2"
$ws.Range("C8").Value = "Gpc4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.25983933333333
$ws.Range("H8").Value = 30.779518
$ws.Range("I8").Value = 0.7373940889775011
$ws.Range("J8").Value = 0.737394088977501
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.04982966666666667
$ws.Range("N8").Value = 0.149489
$ws.Range("O8").Value = 0.0008753708463802955
$ws.Range("P8").Value = 0.0008753708463802955
$ws.Range("Q8").Value = 0.5112443740335556
$ws.Range("R8").Value = 4.601199366302001
$ws.Range("S8").Value = 0.0006454932877840621
$ws.Range("T8").Value = 0.000645493287784062

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fgf2"
$ws.Range("C9").Value = "Gpc4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.25983933333333
$ws.Range("H9").Value = 30.779518
$ws.Range("I9").Value = 0.7373940889775011
$ws.Range("J9").Value = 0.737394088977501
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.888937
$ws.Range("N9").Value = 29.666811
$ws.Range("O9").Value = 0.1737215544586843
$ws.Range("P9").Value = 0.1737215544586843
$ws.Range("Q9").Value = 101.4589047974553
$ws.Range("R9").Value = 913.130143177098
$ws.Range("S9").Value = 0.1281012473858169
$ws.Range("T9").Value = 0.1281012473858169

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Fgf2"
$ws.Range("C10").Value = "Gpc4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.383875
$ws.Range("H10").Value = 1.151625
$ws.Range("I10").Value = 0.02758982345723265
$ws.Range("J10").Value = 0.02758982345723265
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.269799333333333
$ws.Range("N10").Value = 27.809398
$ws.Range("O10").Value = 0.1628450003986012
$ws.Range("P10").Value = 0.1628450003986012
$ws.Range("Q10").Value = 3.558444219083333
$ws.Range("R10").Value = 32.02599797175001
$ws.Range("S10").Value = 0.004492864811890389
$ws.Range("T10").Value = 0.004492864811890389

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Fgf2"
$ws.Range("C11").Value = "Gpc4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.383875
$ws.Range("H11").Value = 1.151625
$ws.Range("I11").Value = 0.02758982345723265
$ws.Range("J11").Value = 0.02758982345723265
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 37.71549866666666
$ws.Range("N11").Value = 113.146496
$ws.Range("O11").Value = 0.6625580742963342
$ws.Range("P11").Value = 0.6625580742963342
$ws.Range("Q11").Value = 14.47803705066667
$ws.Range("R11").Value = 130.302333456
$ws.Range("S11").Value = 0.01827986029999989
$ws.Range("T11").Value = 0.01827986029999989

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Fgf2"
$ws.Range("C12").Value = "Gpc4"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.383875
$ws.Range("H12").Value = 1.151625
$ws.Range("I12").Value = 0.02758982345723265
$ws.Range("J12").Value = 0.02758982345723265
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.04982966666666667
$ws.Range("N12").Value = 0.149489
$ws.Range("O12").Value = 0.0008753708463802955
$ws.Range("P12").Value = 0.0008753708463802955
$ws.Range("Q12").Value = 0.01912836329166667
$ws.Range("R12").Value = 0.172155269625
$ws.Range("S12").Value = 0.00002415132711124068
$ws.Range("T12").Value = 0.00002415132711124068

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Fgf2"
$ws.Range("C13").Value = "Gpc4"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.383875
$ws.Range("H13").Value = 1.151625
$ws.Range("I13").Value = 0.02758982345723265
$ws.Range("J13").Value = 0.02758982345723265
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 9.888937
$ws.Range("N13").Value = 29.666811
$ws.Range("O13").Value = 0.1737215544586843
$ws.Range("P13").Value = 0.1737215544586843
$ws.Range("Q13").Value = 3.796115690875
$ws.Range("R13").Value = 34.16504121787501
$ws.Range("S13").Value = 0.004792947018231128
$ws.Range("T13").Value = 0.004792947018231129

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Fgf2"
$ws.Range("C14").Value = "Gpc4"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.264687
$ws.Range("H14").Value = 6.794061
$ws.Range("I14").Value = 0.162767344880208
$ws.Range("J14").Value = 0.162767344880208
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 9.269799333333333
$ws.Range("N14").Value = 27.809398
$ws.Range("O14").Value = 0.1628450003986012
$ws.Range("P14").Value = 0.1628450003986012
$ws.Range("Q14").Value = 20.99319404280866
$ws.Range("R14").Value = 188.938746385278
$ws.Range("S14").Value = 0.02650584834189674
$ws.Range("T14").Value = 0.02650584834189674

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Fgf2"
$ws.Range("C15").Value = "Gpc4"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.264687
$ws.Range("H15").Value = 6.794061
$ws.Range("I15").Value = 0.162767344880208
$ws.Range("J15").Value = 0.162767344880208
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 37.71549866666666
$ws.Range("N15").Value = 113.146496
$ws.Range("O15").Value = 0.6625580742963342
$ws.Range("P15").Value = 0.6625580742963342
$ws.Range("Q15").Value = 85.41379952891732
$ws.Range("R15").Value = 768.724195760256
$ws.Range("S15").Value = 0.1078428185821579
$ws.Range("T15").Value = 0.1078428185821579

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Fgf2"
$ws.Range("C16").Value = "Gpc4"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.264687
$ws.Range("H16").Value = 6.794061
$ws.Range("I16").Value = 0.162767344880208
$ws.Range("J16").Value = 0.162767344880208
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.04982966666666667
$ws.Range("N16").Value = 0.149489
$ws.Range("O16").Value = 0.0008753708463802955
$ws.Range("P16").Value = 0.0008753708463802955
$ws.Range("Q16").Value = 0.1128485983143333
$ws.Range("R16").Value = 1.015637384829
$ws.Range("S16").Value = 0.0001424817884508611
$ws.Range("T16").Value = 0.0001424817884508611

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Fgf2"
$ws.Range("C17").Value = "Gpc4"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2.264687
$ws.Range("H17").Value = 6.794061
$ws.Range("I17").Value = 0.162767344880208
$ws.Range("J17").Value = 0.162767344880208
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 9.888937
$ws.Range("N17").Value = 29.666811
$ws.Range("O17").Value = 0.1737215544586843
$ws.Range("P17").Value = 0.1737215544586843
$ws.Range("Q17").Value = 22.395347067719
$ws.Range("R17").Value = 201.558123609471
$ws.Range("S17").Value = 0.0282761961677025
$ws.Range("T17").Value = 0.02827619616770251
